$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (94, 95) after the existing data (which ends at row 93).
$ws.Range("A94").Value = 93
$ws.Range("B94").Value = 1
$ws.Range("C94").Value = "2024-06-16 21:13:06"
$ws.Range("D94").Value = 200
$ws.Range("E94").Value = 11

$ws.Range("A95").Value = 94
$ws.Range("B95").Value = 2
$ws.Range("C95").Value = "2024-06-16 21:13:06"
$ws.Range("D95").Value = 200
$ws.Range("E95").Value = 0
